# Adds a "probability" column (E) to the two edge-time sheets
# (test_times_a / test_times_b), with per-edge probability values, and
# updates the sheet selections/active tab to match.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Sheet "test_times_a" (first sheet) - add column E "probability"
# ---------------------------------------------------------------------
$ws1 = $wb.Worksheets.Item("test_times_a")

# Header - copy the header cell style (D1) onto E1 so it keeps the same
# fill as the rest of row 1, then write the label.
$ws1.Range("D1").Copy()
$ws1.Range("E1").PasteSpecial(-4122)   # xlPasteFormats
$excel.CutCopyMode = $false
$ws1.Range("E1").Value = "probability"

$ws1.Range("E2").Value = 0.4
$ws1.Range("E3").Value = 0.3
$ws1.Range("E4").Value = 0.25
$ws1.Range("E5").Value = 0.05

$ws1.Range("E7").Value = 1
$ws1.Range("E8").Value = 1
$ws1.Range("E9").Value = 1
$ws1.Range("E10").Value = 1

$ws1.Range("E12").Value = 1
$ws1.Range("E13").Value = 1
$ws1.Range("E14").Value = 1

$ws1.Range("E16").Value = 0.5
$ws1.Range("E17").Value = 0.3
$ws1.Range("E18").Value = 0.2

$ws1.Range("E20").Value = 0.3
$ws1.Range("E21").Value = 1

$ws1.Range("E23").Value = 0.7
$ws1.Range("E24").Value = 1
$ws1.Range("E25").Value = 1
$ws1.Range("E26").Value = 1

$ws1.Range("E28").Value = 1

# ---------------------------------------------------------------------
# Sheet "test_times_b" (fourth sheet) - add column E "probability"
# ---------------------------------------------------------------------
$ws4 = $wb.Worksheets.Item("test_times_b")

$ws4.Range("D1").Copy()
$ws4.Range("E1").PasteSpecial(-4122)   # xlPasteFormats
$excel.CutCopyMode = $false
$ws4.Range("E1").Value = "probability"

$ws4.Range("E2").Value = 0.4
$ws4.Range("E3").Value = 0.3
$ws4.Range("E4").Value = 0.25
$ws4.Range("E5").Value = 0.05

$ws4.Range("E7").Value = 1
$ws4.Range("E8").Value = 1
$ws4.Range("E9").Value = 1
$ws4.Range("E10").Value = 1

$ws4.Range("E12").Value = 1
$ws4.Range("E13").Value = 1

$ws4.Range("E15").Value = 0.5
$ws4.Range("E16").Value = 0.3
$ws4.Range("E17").Value = 0.2

$ws4.Range("E19").Value = 0.3
$ws4.Range("E20").Value = 1

$ws4.Range("E22").Value = 0.7
$ws4.Range("E23").Value = 1
$ws4.Range("E24").Value = 1
$ws4.Range("E25").Value = 1

$ws4.Range("E27").Value = 1
$ws4.Range("E28").Value = 1
$ws4.Range("E29").Value = 1

$ws4.Range("E31").Value = 1
$ws4.Range("E32").Value = 1
$ws4.Range("E33").Value = 1

$ws4.Range("E35").Value = 1

# ---------------------------------------------------------------------
# View state updates (selection / active sheet)
# ---------------------------------------------------------------------

# "chart_a" keeps its own view, only the selection moves.
$ws2 = $wb.Worksheets.Item("chart_a")
[void]$ws2.Range("A34").Select()

# "test_times_b" is no longer the active tab; only its selection changes.
[void]$ws4.Range("E22").Select()

# "test_times_a" becomes the active/selected tab, selection on E20.
[void]$ws1.Range("E20").Select()
